$wb = $excel.ActiveWorkbook

# --- Config sheet: update values ---
$cfg = $wb.Worksheets.Item("Config")

# Database_Password value: Khanh@2721 -> sa (keep the existing mailto hyperlink,
# but it should now show "Khanh@2721" as the display text override since the
# cell text no longer matches the link target)
$cfg.Range("B18").Value = "sa"
$cfg.Hyperlinks.Item(1).TextToDisplay = "Khanh@2721"
# TextToDisplay assignment resets B18's text back to the link text and can
# duplicate hyperlink entries in this runtime, so re-apply the intended value
# and drop any duplicate hyperlink record it may have produced.
if ($cfg.Hyperlinks.Count -gt 1) {
    $cfg.Hyperlinks.Item($cfg.Hyperlinks.Count).Delete()
}
$cfg.Range("B18").Value = "sa"

# Restore the hyperlink font on B18 (Aptos Narrow -> Arial)
$cfg.Range("B18").Font.Name = "Arial"

# Default_Database_Name value: PE_PRN_Sum25B5_WA -> Library
$cfg.Range("B19").Value = "Library"

# Row 18 no longer needs its custom row height override
$cfg.Rows.Item(18).AutoFit()

# Update view state for the Config sheet
$cfg.Select()
$cfg.Range("B25").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

# --- Run sheet becomes the active tab ---
$run = $wb.Worksheets.Item("Run")
$run.Activate()
